# Auto-generated edit script applying the BRVM daily recommandations update
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Sheet "Recommandations": rows 2-11 (sector indices) - D,E numeric updates ---
$ws1.Range("D2").Value = 1085.24
$ws1.Range("E2").Value = 282.18
$ws1.Range("D3").Value = 1039.45
$ws1.Range("E3").Value = 266.83
$ws1.Range("D4").Value = 761.9
$ws1.Range("E4").Value = 184.63
$ws1.Range("D5").Value = 760
$ws1.Range("E5").Value = 188.36
$ws1.Range("D6").Value = 645.9299999999999
$ws1.Range("E6").Value = 163.5
$ws1.Range("D7").Value = 617.8200000000001
$ws1.Range("E7").Value = 154.66
$ws1.Range("D8").Value = 594.46
$ws1.Range("E8").Value = 150.5
$ws1.Range("D9").Value = 505.71
$ws1.Range("E9").Value = 127.42
$ws1.Range("D10").Value = 476.8
$ws1.Range("E10").Value = 119.42
$ws1.Range("D11").Value = 402.42
$ws1.Range("E11").Value = 101.16

# --- Row 12 (UNILEVER CI): B,D,E numeric updates ---
$ws1.Range("B12").Value = 4
$ws1.Range("D12").Value = 27.32
$ws1.Range("E12").Value = 6.88

# --- Rows 13-32: full re-sort/update of individual stock rows ---
$ws1.Range("A13").Value = 'ECOBANK TRANS. INCORP. TG (ETIT)'
$ws1.Range("B13").Value = 3
$ws1.Range("C13").Value = 0
$ws1.Range("D13").Value = 18.16
$ws1.Range("E13").Value = 6.9
$ws1.Range("F13").Value = '🟢 Achat'
$ws1.Range("G13").Value = '✅ Renforcer'
$ws1.Range("A14").Value = 'CORIS BANK INTERNATIONAL (CBIBF)'
$ws1.Range("B14").Value = 1
$ws1.Range("C14").Value = 0
$ws1.Range("D14").Value = 7.5
$ws1.Range("E14").Value = 7.5
$ws1.Range("F14").Value = '🟡 Observer'
$ws1.Range("G14").Value = '➖ Neutre'
$ws1.Range("A15").Value = 'SMB CI (SMBC)'
$ws1.Range("B15").Value = 1
$ws1.Range("C15").Value = 0
$ws1.Range("D15").Value = 7.49
$ws1.Range("E15").Value = 7.49
$ws1.Range("F15").Value = '🟡 Observer'
$ws1.Range("G15").Value = '➖ Neutre'
$ws1.Range("A16").Value = 'SICOR CI (SICC)'
$ws1.Range("B16").Value = 1
$ws1.Range("C16").Value = 0
$ws1.Range("D16").Value = 7.44
$ws1.Range("E16").Value = 7.44
$ws1.Range("F16").Value = '🟡 Observer'
$ws1.Range("G16").Value = '➖ Neutre'
$ws1.Range("A17").Value = 'SUCRIVOIRE (SCRC)'
$ws1.Range("B17").Value = 1
$ws1.Range("C17").Value = 0
$ws1.Range("D17").Value = 7.34
$ws1.Range("E17").Value = 7.34
$ws1.Range("F17").Value = '🟡 Observer'
$ws1.Range("G17").Value = '➖ Neutre'
$ws1.Range("A18").Value = 'SITAB CI (STBC)'
$ws1.Range("B18").Value = 1
$ws1.Range("C18").Value = 0
$ws1.Range("D18").Value = 5.29
$ws1.Range("E18").Value = 5.29
$ws1.Range("F18").Value = '🟡 Observer'
$ws1.Range("G18").Value = '➖ Neutre'
$ws1.Range("A19").Value = 'UNIWAX CI (UNXC)'
$ws1.Range("B19").Value = 2
$ws1.Range("C19").Value = 1
$ws1.Range("D19").Value = 4.1
$ws1.Range("E19").Value = 4.27
$ws1.Range("F19").Value = '🟡 Observer'
$ws1.Range("G19").Value = '👀 À surveiller'
$ws1.Range("A20").Value = 'SOLIBRA CI (SLBC)'
$ws1.Range("B20").Value = 1
$ws1.Range("C20").Value = 0
$ws1.Range("D20").Value = 2.17
$ws1.Range("E20").Value = 2.17
$ws1.Range("F20").Value = '🟡 Observer'
$ws1.Range("G20").Value = '➖ Neutre'
$ws1.Range("A21").Value = 'LOTERIE NATIONALE DU BENIN (LNBB)'
$ws1.Range("B21").Value = 1
$ws1.Range("C21").Value = 0
$ws1.Range("D21").Value = 1.65
$ws1.Range("E21").Value = 1.65
$ws1.Range("F21").Value = '🟡 Observer'
$ws1.Range("G21").Value = '➖ Neutre'
$ws1.Range("A22").Value = 'AFRICA GLOBAL LOGISTICS CI (SDSC)'
$ws1.Range("B22").Value = 0
$ws1.Range("C22").Value = 1
$ws1.Range("D22").Value = -2.29
$ws1.Range("E22").Value = -2.29
$ws1.Range("F22").Value = '🟡 Observer'
$ws1.Range("G22").Value = '➖ Neutre'
$ws1.Range("A23").Value = 'SETAO CI (STAC)'
$ws1.Range("B23").Value = 1
$ws1.Range("C23").Value = 2
$ws1.Range("D23").Value = -2.49
$ws1.Range("E23").Value = -2.99
$ws1.Range("F23").Value = '🟡 Observer'
$ws1.Range("G23").Value = '👀 À surveiller'
$ws1.Range("A24").Value = 'BANK OF AFRICA BF (BOABF)'
$ws1.Range("B24").Value = 0
$ws1.Range("C24").Value = 1
$ws1.Range("D24").Value = -2.72
$ws1.Range("E24").Value = -2.72
$ws1.Range("F24").Value = '🟡 Observer'
$ws1.Range("G24").Value = '➖ Neutre'
$ws1.Range("A25").Value = 'ONATEL BF (ONTBF)'
$ws1.Range("B25").Value = 0
$ws1.Range("C25").Value = 1
$ws1.Range("D25").Value = -2.95
$ws1.Range("E25").Value = -2.95
$ws1.Range("F25").Value = '🟡 Observer'
$ws1.Range("G25").Value = '➖ Neutre'
$ws1.Range("A26").Value = 'TRACTAFRIC MOTORS CI (PRSC)'
$ws1.Range("B26").Value = 0
$ws1.Range("C26").Value = 1
$ws1.Range("D26").Value = -4.67
$ws1.Range("E26").Value = -4.67
$ws1.Range("F26").Value = '🟡 Observer'
$ws1.Range("G26").Value = '➖ Neutre'
$ws1.Range("A27").Value = 'ECOBANK COTE D''''IVOIRE (ECOC)'
$ws1.Range("B27").Value = 1
$ws1.Range("C27").Value = 2
$ws1.Range("D27").Value = -6.45
$ws1.Range("E27").Value = -5.33
$ws1.Range("F27").Value = '🟡 Observer'
$ws1.Range("G27").Value = '👀 À surveiller'
$ws1.Range("A28").Value = 'BERNABE CI (BNBC)'
$ws1.Range("B28").Value = 0
$ws1.Range("C28").Value = 2
$ws1.Range("D28").Value = -7.07
$ws1.Range("E28").Value = -4.13
$ws1.Range("F28").Value = '🟡 Observer'
$ws1.Range("G28").Value = '➖ Neutre'
$ws1.Range("A29").Value = 'FILTISAC CI (FTSC)'
$ws1.Range("B29").Value = 0
$ws1.Range("C29").Value = 1
$ws1.Range("D29").Value = -7.26
$ws1.Range("E29").Value = -7.26
$ws1.Range("F29").Value = '🟡 Observer'
$ws1.Range("G29").Value = '➖ Neutre'
$ws1.Range("A30").Value = 'ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)'
$ws1.Range("B30").Value = 1
$ws1.Range("C30").Value = 2
$ws1.Range("D30").Value = -7.82
$ws1.Range("E30").Value = -7.39
$ws1.Range("F30").Value = '🟡 Observer'
$ws1.Range("G30").Value = '👀 À surveiller'
$ws1.Range("A31").Value = 'EVIOSYS PACKAGING SIEM CI (SEMC)'
$ws1.Range("B31").Value = 1
$ws1.Range("C31").Value = 3
$ws1.Range("D31").Value = -14.72
$ws1.Range("E31").Value = -7.3
$ws1.Range("F31").Value = '🔴 Vente'
$ws1.Range("G31").Value = '⚠️ Risque de décrochage'
$ws1.Range("A32").Value = 'SICABLE CI (CABC)'
$ws1.Range("B32").Value = 0
$ws1.Range("C32").Value = 3
$ws1.Range("D32").Value = -22.38
$ws1.Range("E32").Value = -7.45
$ws1.Range("F32").Value = '🔴 Vente'
$ws1.Range("G32").Value = '⚠️ Risque de décrochage'

# --- Remove now-obsolete rows 33-34 (table shrank from 34 to 32 data-company rows) ---
$ws1.Rows("33:34").Delete()

# --- Sheet "Top_YTD": B2:B11 numeric updates ---
$ws2.Range("B2").Value = 18892.46
$ws2.Range("B3").Value = 16665.21
$ws2.Range("B4").Value = 7013.43
$ws2.Range("B5").Value = 6972.29
$ws2.Range("B6").Value = 4574.58
$ws2.Range("B7").Value = 4092.18
$ws2.Range("B8").Value = 3720.17
$ws2.Range("B9").Value = 2528.46
$ws2.Range("B10").Value = 2208.61
$ws2.Range("B11").Value = 1519.42
